# Adds two new columns, I ("I0") and J ("IF"), to the sheet, filling
# headers in row 1 and numeric values for rows 2-64 (data rows 0-62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (I1, J1) -------------------------------------------------
# Copy the formatting of the existing header cell H1 (bold font, thin box
# border, centered alignment) onto the two new header cells, then set their
# text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells (I2:J64) ----------------------------------------------------
$values = @(
        @(8,8),
        @(7,8),
        @(8,8),
        @(8,8),
        @(7,8),
        @(8,8),
        @(7,8),
        @(8,8),
        @(8,8),
        @(8,8),
        @(7,8),
        @(9,9),
        @(7,7),
        @(8,8),
        @(8,8),
        @(8,8),
        @(8,8),
        @(9,9),
        @(10,10),
        @(9,9),
        @(9,10),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(9,9),
        @(7,7),
        @(8,9),
        @(8,8),
        @(9,9),
        @(7,8),
        @(9,9),
        @(8,8),
        @(9,9),
        @(7,8),
        @(9,9),
        @(9,9),
        @(8,8),
        @(8,8),
        @(7,8),
        @(8,8),
        @(7,8),
        @(8,8),
        @(6,7),
        @(6,7),
        @(6,7),
        @(6,6),
        @(7,7),
        @(8,8),
        @(6,6),
        @(8,8),
        @(8,8),
        @(9,9),
        @(7,7),
        @(7,7),
        @(7,7),
        @(5,5)
    )

$startRow = 2
for ($idx = 0; $idx -lt $values.Count; $idx++) {
    $row = $startRow + $idx
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
